$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: About
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A10").Value = "We do not apply this potential to non-manufacturing indutries due to their lack of large, heat-generating,"

# ---------------------------------------------------------------------------
# Sheet: Data
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("A7").Value = "Percentage Savings in 2030 (model end year)"

# ---------------------------------------------------------------------------
# Sheet: PPRiFUfICaWHR
# ---------------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("PPRiFUfICaWHR")

# Remove the old 9-row layout and rebuild it (also drops the stale per-row
# height override of 14.45 that Excel had stamped on rows 1-9).
$wsP.Range("A1:B9").EntireRow.Delete()

$wsP.Range("A1").Value = "Unit: dimensionless"
$wsP.Range("A1").Font.Italic = $true
$wsP.Range("B1").Value = "Pot Perc Red in Fuel Use"
$wsP.Range("B1").HorizontalAlignment = -4152
$wsP.Range("B1").WrapText = $true

$industries = @(
  "agriculture and forestry 01T03",
  "coal mining 05",
  "oil and gas extraction 06",
  "other mining and quarrying 07T08",
  "food beverage and tobacco 10T12",
  "textiles apparel and leather 13T15",
  "wood products 16",
  "pulp paper and printing 17T18",
  "refined petroleum and coke 19",
  "chemicals 20",
  "rubber and plastic products 22",
  "glass and glass products 231",
  "cement and other nonmetallic minerals 239",
  "iron and steel 241",
  "other metals 242",
  "metal products except machinery and vehicles 25",
  "computers and electronics 26",
  "appliances and electrical equipment 27",
  "other machinery 28",
  "road vehicles 29",
  "nonroad vehicles 30",
  "other manufacturing 31T33",
  "energy pipelines and gas processing 352T353",
  "water and waste 36T39",
  "construction 41T43"
)

# Rows 2-5 and 24-26 (relative to the new table, i.e. industries 1-4 and
# 23-25) are hard-coded to 0; the rest pull the model's Data!B$7 figure.
$zeroRows = @(1,2,3,4,23,24,25)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $wsP.Range("A$row").Value = $industries[$i]
    if ($zeroRows -contains ($i + 1)) {
        $wsP.Range("B$row").Value = 0
    } else {
        $wsP.Range("B$row").Formula = "=Data!B`$7"
    }
}

$wsP.Columns.Item(1).ColumnWidth = 47

# ---------------------------------------------------------------------------
# Tab / selection state
# ---------------------------------------------------------------------------
$wsAbout.Activate()
